{"js": "// Add the GitHub repository link as a new run in the last (currently\n// empty) bullet paragraph of the document body \u2014 the sub-bullet under\n// \"My GitHub Repository link to Module 7:\" \u2014 matching the run\n// formatting (Times New Roman, 12pt / sz 24) already used for the rest\n// of that list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The target is the final paragraph in the document body: an empty\n// \"ListParagraph\" (ilvl 1) bullet immediately before the closing\n// sectPr, meant to hold the repo URL.\nconst items = paragraphs.items;\nconst target = items[items.length - 1];\n\nconst url = \"https://github.com/RUKtheCROOK/CSD-310/tree/main/Module%207\";\n\n// Insert the URL text into the (empty) paragraph and pick up the\n// resulting range so we can stamp the run formatting onto it.\nconst range = target.insertText(url, Word.InsertLocation.replace);\n\n// Match the rPr used elsewhere in this list: Times New Roman for\n// ascii/hAnsi/cs, 12pt (half-points 24) for both sz and szCs.\nrange.font.name = \"Times New Roman\";\nrange.font.nameAscii = \"Times New Roman\";\nrange.font.nameBidirectional = \"Times New Roman\";\nrange.font.size = 12;\nrange.font.sizeBidirectional = 12;\n\nawait context.sync();\n", "ps1": "# Add the GitHub repository link as a new run in the last (currently\n# empty) bullet paragraph of the document body -- the sub-bullet that\n# follows \"My GitHub Repository link to Module 7:\" -- matching the run\n# formatting (Times New Roman, 12pt / sz 24) already used for the rest\n# of that list.\n\n$d = $word.ActiveDocument\n$url = \"https://github.com/RUKtheCROOK/CSD-310/tree/main/Module%207\"\n\n# Locate the \"My GitHub Repository link to Module 7:\" paragraph, then\n# target the paragraph right after it -- the empty sub-bullet meant to\n# hold the URL. Fall back to the very last paragraph in the document\n# (same empty sub-bullet) if the label text can't be found for some\n# reason, so the edit is still applied.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*My GitHub Repository link*\" -and $i -lt $count) {\n        $target = $d.Paragraphs.Item($i + 1)\n    }\n}\nif ($target -eq $null) {\n    $target = $d.Paragraphs.Last\n}\n\n$rng = $target.Range\n$rng.Text = $url\n\n# Match the rPr used elsewhere in this list: Times New Roman for\n# ascii/hAnsi/cs, 12pt (half-points 24) for both sz and szCs.\n$rng.Font.Name = \"Times New Roman\"\n$rng.Font.NameAscii = \"Times New Roman\"\n$rng.Font.NameBi = \"Times New Roman\"\n$rng.Font.Size = 12\n$rng.Font.SizeBi = 12\n"}
